{"js": "// Update the 25 \"three-digit \u00f7 one-digit\" answer cells to the new\n// generated values. Each old answer string is unique in the document,\n// so a targeted body.search() + insertText(..., replace) per pair is\n// safe and only rewrites the matched run's text.\n\nconst pairs = [\n  [\"387\u00f73=129, 0\", \"454\u00f77=64, 6\"],\n  [\"880\u00f77=125, 5\", \"108\u00f73=36, 0\"],\n  [\"825\u00f73=275, 0\", \"858\u00f72=429, 0\"],\n  [\"505\u00f72=252, 1\", \"394\u00f79=43, 7\"],\n  [\"340\u00f75=68, 0\", \"975\u00f72=487, 1\"],\n  [\"516\u00f76=86, 0\", \"598\u00f73=199, 1\"],\n  [\"783\u00f75=156, 3\", \"544\u00f79=60, 4\"],\n  [\"218\u00f76=36, 2\", \"395\u00f79=43, 8\"],\n  [\"237\u00f72=118, 1\", \"919\u00f76=153, 1\"],\n  [\"855\u00f78=106, 7\", \"763\u00f73=254, 1\"],\n  [\"779\u00f78=97, 3\", \"188\u00f73=62, 2\"],\n  [\"850\u00f73=283, 1\", \"219\u00f73=73, 0\"],\n  [\"934\u00f75=186, 4\", \"990\u00f79=110, 0\"],\n  [\"966\u00f72=483, 0\", \"523\u00f73=174, 1\"],\n  [\"672\u00f73=224, 0\", \"631\u00f77=90, 1\"],\n  [\"931\u00f74=232, 3\", \"981\u00f73=327, 0\"],\n  [\"744\u00f78=93, 0\", \"455\u00f79=50, 5\"],\n  [\"782\u00f75=156, 2\", \"311\u00f76=51, 5\"],\n  [\"364\u00f74=91, 0\", \"555\u00f79=61, 6\"],\n  [\"848\u00f72=424, 0\", \"114\u00f75=22, 4\"],\n  [\"689\u00f78=86, 1\", \"775\u00f77=110, 5\"],\n  [\"810\u00f75=162, 0\", \"475\u00f73=158, 1\"],\n  [\"516\u00f74=129, 0\", \"857\u00f76=142, 5\"],\n  [\"566\u00f72=283, 0\", \"322\u00f78=40, 2\"],\n  [\"223\u00f72=111, 1\", \"980\u00f76=163, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit \u00f7 one-digit\" answer cells to the new\n# generated values. Each old answer string is unique in the document,\n# so a single targeted Find/Replace (wdReplaceOne) per pair is safe and\n# avoids touching any other text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"387\u00f73=129, 0\", \"454\u00f77=64, 6\"),\n  @(\"880\u00f77=125, 5\", \"108\u00f73=36, 0\"),\n  @(\"825\u00f73=275, 0\", \"858\u00f72=429, 0\"),\n  @(\"505\u00f72=252, 1\", \"394\u00f79=43, 7\"),\n  @(\"340\u00f75=68, 0\", \"975\u00f72=487, 1\"),\n  @(\"516\u00f76=86, 0\", \"598\u00f73=199, 1\"),\n  @(\"783\u00f75=156, 3\", \"544\u00f79=60, 4\"),\n  @(\"218\u00f76=36, 2\", \"395\u00f79=43, 8\"),\n  @(\"237\u00f72=118, 1\", \"919\u00f76=153, 1\"),\n  @(\"855\u00f78=106, 7\", \"763\u00f73=254, 1\"),\n  @(\"779\u00f78=97, 3\", \"188\u00f73=62, 2\"),\n  @(\"850\u00f73=283, 1\", \"219\u00f73=73, 0\"),\n  @(\"934\u00f75=186, 4\", \"990\u00f79=110, 0\"),\n  @(\"966\u00f72=483, 0\", \"523\u00f73=174, 1\"),\n  @(\"672\u00f73=224, 0\", \"631\u00f77=90, 1\"),\n  @(\"931\u00f74=232, 3\", \"981\u00f73=327, 0\"),\n  @(\"744\u00f78=93, 0\", \"455\u00f79=50, 5\"),\n  @(\"782\u00f75=156, 2\", \"311\u00f76=51, 5\"),\n  @(\"364\u00f74=91, 0\", \"555\u00f79=61, 6\"),\n  @(\"848\u00f72=424, 0\", \"114\u00f75=22, 4\"),\n  @(\"689\u00f78=86, 1\", \"775\u00f77=110, 5\"),\n  @(\"810\u00f75=162, 0\", \"475\u00f73=158, 1\"),\n  @(\"516\u00f74=129, 0\", \"857\u00f76=142, 5\"),\n  @(\"566\u00f72=283, 0\", \"322\u00f78=40, 2\"),\n  @(\"223\u00f72=111, 1\", \"980\u00f76=163, 2\")\n)\n\nforeach ($p in $pairs) {\n  $rng = $d.Content\n  $rng.Find.Text = $p[0]\n  $rng.Find.Execute($p[0], $false, $false, $false, $false, $false, $true, 1, $false, $p[1], 1)\n}\n"}
